# Added newly tested devices
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Fix manufacturer typo on existing row (Motorala -> Motorola)
$ws.Range("A7").Value = "Motorola"

# Append the newly tested Huawei device as row 8
$ws.Range("A8").Value = "Huawei"
$ws.Range("B8").Value = "ALE-L21"
$ws.Range("C8").Value = "P8 Lite"
$ws.Range("D8").Value = "Android 6"
$ws.Range("E8").Value = "WIFI MAC not possible`nSerial Number not possible"
$ws.Range("E8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 30

# Turn on the AutoFilter for the header row
$ws.Range("A4:E4").AutoFilter() | Out-Null
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Tabelle1!`$A`$4:`$E`$4")
$filterName.Visible = $false
